$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '305.05'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '0.35%'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '35.62'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '-4.11%'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.055'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '0.45%'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.07882'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '0.30%'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '2.137'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '-3.22%'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '7.941'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '-0.59%'
$ws.Range('B8').Value = 'GateToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '4.131'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '2.95%'
$ws.Range('B9').Value = 'MXToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.9240'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '-0.30%'
$ws.Range('B10').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C10').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.09689'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '-1.60%'
$ws.Range('B11').Value = 'WazirX'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.1839'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '-2.13%'
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08674'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '0.67%'
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.03570'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '-3.44%'
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.09920'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '0.06%'
$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.001439'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '-2.08%'
$ws.Range('B16').Value = 'CoinExToken'
$ws.Range('C16').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.04563'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '-1.33%'
$ws.Range('B17').Value = 'TigerCash'
$ws.Range('C17').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.005658'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '-0.50%'
$ws.Range('B18').Value = 'LEO'
$ws.Range('C18').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.474'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '0.18%'
$ws.Range('B19').Value = 'BTSEToken'
$ws.Range('C19').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '2.753'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '22.23%'
$ws.Range('B20').Value = 'BitpandaEcosystemToken'
$ws.Range('C20').Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.3370'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '-1.15%'
$ws.Range('B21').Value = 'ProBitToken'
$ws.Range('C21').Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.1346'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '3.52%'
$ws.Range('B22').Value = 'MCDex'
$ws.Range('C22').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.156'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '8.11%'
$ws.Range('B23').Value = 'ZBToken'
$ws.Range('C23').Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.2210'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '0.51%'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '-1.82%'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.004831'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '7.69%'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0001302'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '-7.33%'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0004755'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '74.92%'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01855'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '1.16%'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.04739'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '-0.39%'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.007787'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '-2.61%'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1389'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '-0.55%'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.007754'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '2.10%'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.002166'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '2.61%'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.01133'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '12.18%'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00006322'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '0.65%'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.00000000751'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '0.09%'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '0.12%'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '51.28'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '66.98%'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.001902'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '-29.28%'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.00002103'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '0.09%'

Write-Output "Applied 108 cell updates."
